{"js": "// Renumber the \"Online Supplement\" figure captions from the old\n// \"O<n>\" scheme to the new \"S<n+8>\" scheme (O1->S9, O2->S10, ... O11->S19).\n//\n// The captions appear either as a bare \"Figure O<n>:\" lead-in inside a\n// single run, or (for a few figures) already split across multiple runs\n// such as \"Figure \" / \"O<n>\" / \": <rest>\". In every case the only text\n// that actually changes is the \"O<n>\" figure-number token itself, so we\n// locate that token paragraph-by-paragraph (to avoid \"O1\" accidentally\n// matching inside \"O10\"/\"O11\") and replace just that token, leaving the\n// surrounding run formatting (bold, Times New Roman) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Matches a figure-number token \"O\" followed by digits, either right\n// after \"Figure \" or standing completely alone in its own paragraph\n// (the latter covers captions that already had \"Figure \"/\"O#\"/\": ...\"\n// split into separate runs before this edit).\nconst leadPattern = /^Figure (O(\\d+))\\b/;\nconst bareTokenPattern = /^(O(\\d+))$/;\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  const leadMatch = text.match(leadPattern);\n  const bareMatch = text.match(bareTokenPattern);\n  const match = leadMatch || bareMatch;\n  if (!match) continue;\n\n  const oldToken = match[1];         // e.g. \"O1\"\n  const figureNum = parseInt(match[2], 10); // e.g. 1\n  const newToken = \"S\" + (figureNum + 8);   // e.g. \"S9\"\n\n  const paraRange = paragraph.getRange();\n  const found = paraRange.search(oldToken, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    found.items[0].insertText(newToken, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Renumber the \"Online Supplement\" figure captions from the old\n# \"O<n>\" scheme to the new \"S<n+8>\" scheme (O1->S9, O2->S10, ... O11->S19).\n#\n# The captions appear either as a \"Figure O<n>:\" lead-in inside a single\n# run, or (for a few figures) already split across multiple runs such as\n# \"Figure \" / \"O<n>\" / \": <rest>\". In every case the only text that\n# actually changes is the \"O<n>\" figure-number token itself, so we find\n# that token paragraph-by-paragraph (to avoid \"O1\" accidentally matching\n# inside \"O10\"/\"O11\") and overwrite just that token via a precise\n# character-offset Range, leaving the rest of the paragraph untouched.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n\n    $oldToken = $null\n    $figureNum = 0\n\n    if ($text -match \"^Figure (O(\\d+))\\b\") {\n        $oldToken = $matches[1]\n        $figureNum = [int]$matches[2]\n    } elseif ($text -match \"^(O(\\d+))\\r?$\") {\n        $oldToken = $matches[1]\n        $figureNum = [int]$matches[2]\n    }\n\n    if ($oldToken) {\n        $newToken = \"S\" + ($figureNum + 8)\n\n        $offset = $text.IndexOf($oldToken)\n        $start = $p.Range.Start + $offset\n        $end = $start + $oldToken.Length\n\n        $tokenRange = $d.Range($start, $end)\n        $tokenRange.Text = $newToken\n    }\n}\n"}
